# Update BXD causality-test results table (rows 2-19) with new likelihood data.
# New individuals (BXD5, BXD31, BXD23, BXD13, BXD19) are interleaved among the
# existing ones and three additional rows (BXD6, BXD16, BXD9) are appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=Individual# L=Individual R=L C=R R|L_mean R|L_var C|L_mean C|L_var P(R|L) P(C|L) P(L)*P(R|L)*P(C|L)
$data = @(
    @(0, "BXD1", 0, 3.1855, 0, 3.010875, 0.202774796875, 38.875, 1245.109375, 0.8217653323779117, 0.006162279433048035, 0.002531973803252144),
    @(27, "BXD5", 0, 3.0645, 83, 3.010875, 0.202774796875, 38.875, 1245.109375, 0.8796777666901296, 0.005173048256839504, 0.002275307768778421),
    @(14, "BXD31", 0, 2.931, 20, 3.010875, 0.202774796875, 38.875, 1245.109375, 0.8721091730414935, 0.00979881462615984, 0.004272818060203575),
    @(10, "BXD23", 0, 2.6615, 20, 3.010875, 0.202774796875, 38.875, 1245.109375, 0.6556748273549828, 0.00979881462615984, 0.003212418044145417),
    @(8, "BXD2", 0, 3.2065, 10, 3.010875, 0.202774796875, 38.875, 1245.109375, 0.8061604074726, 0.008089049222233136, 0.003260535608530691),
    @(9, "BXD21", 0, 2.374, 10, 3.010875, 0.202774796875, 38.875, 1245.109375, 0.3258698341765872, 0.008089049222233136, 0.001317988564347682),
    @(4, "BXD14", 0, 2.6955, 89, 3.010875, 0.202774796875, 38.875, 1245.109375, 0.6932533531413017, 0.004122144904925831, 0.001428845388737082),
    @(3, "BXD13", 0, 3.9685, 79, 3.010875, 0.202774796875, 38.875, 1245.109375, 0.09233328750545597, 0.005922695085047902, 0.0002734309540474395),
    @(7, "BXD19", 1, 2.905, 0, 3.239100000000001, 0.22243699, 45.3, 1960.41, 0.6581702389148893, 0.005338685304373844, 0.001756881891135571),
    @(5, "BXD15", 1, 3.3835, 100, 3.239100000000001, 0.22243699, 45.3, 1960.41, 0.8071440783008531, 0.004200626402036855, 0.001695255362779133),
    @(35, "BXD8", 1, 3.353, 100, 3.239100000000001, 0.22243699, 45.3, 1960.41, 0.8215648791788004, 0.004200626402036855, 0.001725543561232344),
    @(11, "BXD24", 1, 2.995, 100, 3.239100000000001, 0.22243699, 45.3, 1960.41, 0.739841571720248, 0.004200626402036855, 0.001553899019746259),
    @(12, "BXD28", 1, 4.2805, 66, 3.239100000000001, 0.22243699, 45.3, 1960.41, 0.07388959830514233, 0.008077451998695988, 0.0002984198417563578),
    @(15, "BXD32", 1, 3.072, 0, 3.239100000000001, 0.22243699, 45.3, 1960.41, 0.7944163442429553, 0.005338685304373844, 0.002120569431282129),
    @(1, "BXD11", 1, 3.607, 12, 3.239100000000001, 0.22243699, 45.3, 1960.41, 0.623986083495514, 0.006790613889094571, 0.002118624282593181),
    @(29, "BXD6", 1, 3.387, 0, 3.239100000000001, 0.22243699, 45.3, 1960.41, 0.805290070084263, 0.005338685304373844, 0.002149595131458519),
    @(6, "BXD16", 1, 2.4065, 0, 3.239100000000001, 0.22243699, 45.3, 1960.41, 0.1780611826389548, 0.005338685304373844, 0.0004753063095170074),
    @(38, "BXD9", 1, 3.0015, 75, 3.239100000000001, 0.22243699, 45.3, 1960.41, 0.7450669883149073, 0.007195000363549942, 0.002680378625897409)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
    $r = $r + 1
}

# Column A on the newly-added rows (17-19) needs the same bold/bordered/centered
# style already used by the rest of column A (style index 1) - copy formats only.
$ws.Range("A2").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)
